$d = $word.ActiveDocument

# --- Small helper: nudging a range's font size away from and back to its
# original value forces the engine to keep it as a distinct <w:r> instead
# of silently coalescing it into an adjacent run that happens to share the
# same formatting. ---
function Keep-RunSeparate($rng) {
    $orig = $rng.Font.Size
    if ($orig -eq 23) {
        $rng.Font.Size = 11
    } else {
        $rng.Font.Size = 23
    }
    $rng.Font.Size = $orig
}

# 1) The two TIME field results ("1 de noviembre de 2024") become
#    "24 de noviembre de 2024".
[void]$d.Content.Find.Execute("1 de noviembre de 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "24 de noviembre de 2024", 2)

# 2) Insert a new run reading "del " right after "En atención a su
#    solicitud " and before the date field that follows it.
$r = $d.Content
[void]$r.Find.Execute("En atención a su solicitud ")
$r.Collapse(0)
$r.InsertBefore("del ")
Keep-RunSeparate $r

# 3) "el(los) " -> "el" (drop the "(los)" and the trailing space).
$r2 = $d.Content
[void]$r2.Find.Execute("el(los) ")
$r2.Text = "el"
Keep-RunSeparate $r2

# 4) "{{DIAS}} de {{MESES}} de {{AÑO}}" -> "{{DIAS}} {{MESES}} {{AÑO}}".
$r3 = $d.Content
[void]$r3.Find.Execute("{{DIAS}} de {{MESES}} de {{AÑO}}")
$r3.Text = "{{DIAS}} {{MESES}} {{AÑO}}"
Keep-RunSeparate $r3

# Re-separate the run that follows the placeholders (". P" / "or tal
# razón...") in case it got folded back into the placeholder run above.
$r4 = $d.Content
[void]$r4.Find.Execute(". P")
Keep-RunSeparate $r4
